# "Tambah Peserta" - add six new MENDAFTAR (selection) participants to the
# pesertaseleksidb workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # MENDAFTAR
$ws2 = $wb.Worksheets.Item(2)   # PANAT2014

# New participants: NIM, Nama, CP (phone), Tanggal Seleksi I (date serial)
$peserta = @(
  @("14111-3473", "Ami Aivya Sibarani",     "085297462566", 42187),
  @("14111-3619", "Eben Roy H. Silalahi",   "082165152615", 42187),
  @("14811-1411", "Dian Kristian Silaban",  "085206514334", 42188),
  @("14811-1497", "Widia Marito Manulang",  "082276099679", 42188),
  @("14811-1705", "Shinta Warni Meliala",   "085761761106", 42188),
  @("14811-0459", "Yuliana Siahaan",        "085372781223", 42188)
)

$row = 11
foreach ($p in $peserta) {
  $ws1.Cells.Item($row, 1).Value = $p[0]   # A: NIM
  $ws1.Cells.Item($row, 2).Value = $p[1]   # B: Nama
  $ws1.Cells.Item($row, 4).Value = $p[2]   # D: CP
  $ws1.Cells.Item($row, 5).Value = $p[3]   # E: Tanggal Seleksi I
  $row++
}

# Update the lingering selection on the PANAT2014 sheet's frozen bottom pane,
# then restore MENDAFTAR as the active sheet with its own cursor moved to A17
# (just past the newly-entered rows).
$ws2.Range("B26").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A17").Select() | Out-Null
